# Updated symbol list on Wed Dec 28 23:28:10 UTC 2022 with GitHub Actions
#
# Refresh the "Price" (column D) and "Volume(1h)" (column E) values that
# changed for this crawl. Price values look numeric, so they must be
# written as text (leading apostrophe) to preserve the original inline
# string cell type, then the style is reset back to Normal so no stray
# "quote prefix" / text number-format gets attached to the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Cell,
        [string]$Text
    )
    $ws.Range($Cell).Value = "'" + $Text
    $ws.Range($Cell).Style = "Normal"
}

# Column D (Price) updates
Set-TextValue "D2"  "243.22"
Set-TextValue "D3"  "23.80"
Set-TextValue "D4"  "5.232"
Set-TextValue "D5"  "0.05757"
Set-TextValue "D6"  "6.402"
Set-TextValue "D7"  "3.229"
Set-TextValue "D8"  "0.8123"
Set-TextValue "D9"  "0.8866"
Set-TextValue "D10" "0.1373"
Set-TextValue "D11" "0.07074"
Set-TextValue "D13" "0.03044"
Set-TextValue "D14" "0.09319"
Set-TextValue "D15" "3.810"
Set-TextValue "D16" "0.001515"
Set-TextValue "D18" "0.0006004"
Set-TextValue "D19" "0.006189"
Set-TextValue "D20" "0.001238"
Set-TextValue "D23" "3.547"
Set-TextValue "D24" "2.147"
Set-TextValue "D25" "0.3165"
Set-TextValue "D40" "0.03728"
Set-TextValue "D41" "0.006264"
Set-TextValue "D42" "0.1046"
Set-TextValue "D43" "0.002408"
Set-TextValue "D44" "0.007121"
Set-TextValue "D45" "0.00005297"
Set-TextValue "D48" "0.002333"

# Column E (Volume(1h)) updates
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("E47").Value = "46CoinbaseStockTokenCOIN"
$ws.Range("E48").Value = "47BOLOBOLOBestin24h"
